$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.804.43"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "2.254.02"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.33"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.82"
$ws.Range("E7").Value = "  +5.40%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.652"
$ws.Range("E9").Value = "  +14.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.73"
$ws.Range("E10").Value = "  +9.04%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.62"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0968"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.41"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "2.591.64"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.77"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "2.243.73"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "42.680.44"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.22"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.20"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.67"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.71"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.04"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.130"
$ws.Range("E33").Value = "  +10.69%  "
$ws.Range("E34").Value = "  +12.69%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.126"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.68"
$ws.Range("E37").Value = "  +7.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.14"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0320"
$ws.Range("E40").Value = "  +7.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.29"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.49"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.84"
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.02"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.01"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.01"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("E51").Value = "  +1.99%  "
